# Fix a bug in the trading module: a trade that was missing from the
# exported trade log. Append it as a new row (row 5) to the trades sheet,
# mirroring the layout of the existing trade rows.
#
# Columns: A Date | B Profitable | C Principle | D Start Principle |
#          E BuyPrice | F SellPrice | G IsShortSell | H Price Change % |
#          I Strong trade

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed row 5 from row 4 first so the new row inherits the same cell
# formatting/styles (e.g. the date format on column A and the boolean
# style on column G) instead of creating brand-new style entries.
$ws.Range("A4:I4").Copy($ws.Range("A5:I5"))

# Now overwrite the copied values with the actual trade data for this row.
$ws.Range("A5").Value = 42636.606736111113
$ws.Range("B5").Value = $false
$ws.Range("C5").Value = 9956.0499999999993
$ws.Range("D5").Value = 9974
$ws.Range("E5").Value = 19.29
$ws.Range("F5").Value = 19.22
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = -0.36
$ws.Range("I5").Value = $false
